# Auto-generated script applying cached-value updates per the target diff.
# Each worksheet (by name) has specific cells whose cached numeric values are updated
# to reflect refreshed market data computed by the scheduled runner.
$wb = $excel.ActiveWorkbook

# ---- Worksheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6462.6113
$ws.Range("I40").Value = 3943.1667
$ws.Range("J40").Value = 7722.3335
$ws.Range("K40").Value = 3943.1667
$ws.Range("L40").Value = 7722.3335
$ws.Range("M40").Value = -3768.1667
$ws.Range("N40").Value = -8072.3335
$ws.Range("H64").Value = 8363.637000000001
$ws.Range("I64").Value = 4000
$ws.Range("J64").Value = 12000
$ws.Range("K64").Value = 4000
$ws.Range("L64").Value = 12000
$ws.Range("M64").Value = -3752
$ws.Range("N64").Value = -12496
$ws.Range("H67").Value = 8363.637000000001
$ws.Range("I67").Value = 4000
$ws.Range("J67").Value = 12000
$ws.Range("K67").Value = 4000
$ws.Range("L67").Value = 12000
$ws.Range("M67").Value = -3142
$ws.Range("N67").Value = -13716
$ws.Range("H106").Value = 4737.107
$ws.Range("I106").Value = 4519.087
$ws.Range("K106").Value = 4519.087
$ws.Range("M106").Value = -3888.087
$ws.Range("H137").Value = 34605
$ws.Range("I137").Value = 45854.332
$ws.Range("J137").Value = 6481.6665
$ws.Range("K137").Value = 137562.996
$ws.Range("L137").Value = 19444.9995
$ws.Range("M137").Value = -135012.996
$ws.Range("N137").Value = -24544.9995
$ws.Range("H138").Value = 5628.75
$ws.Range("J138").Value = 6205.1665
$ws.Range("L138").Value = 18615.4995
$ws.Range("N138").Value = -28895.4995

# ---- Worksheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4264.6304
$ws.Range("I32").Value = 3772.7097
$ws.Range("J32").Value = 7037.273
$ws.Range("K32").Value = 3772.7097
$ws.Range("L32").Value = 7037.273
$ws.Range("M32").Value = -3485.7097
$ws.Range("N32").Value = -7611.273
$ws.Range("H63").Value = 8295
$ws.Range("I63").Value = 8886
$ws.Range("J63").Value = 7999.5
$ws.Range("K63").Value = 8886
$ws.Range("L63").Value = 7999.5
$ws.Range("M63").Value = -8200
$ws.Range("N63").Value = -9371.5
$ws.Range("H66").Value = 8295
$ws.Range("I66").Value = 8886
$ws.Range("J66").Value = 7999.5
$ws.Range("K66").Value = 44430
$ws.Range("L66").Value = 39997.5
$ws.Range("M66").Value = -40998
$ws.Range("N66").Value = -46861.5
$ws.Range("H74").Value = 39268.21
$ws.Range("I74").Value = 3961.8572
$ws.Range("J74").Value = 138126
$ws.Range("K74").Value = 3961.8572
$ws.Range("L74").Value = 138126
$ws.Range("M74").Value = -3087.8572
$ws.Range("N74").Value = -139874
$ws.Range("H77").Value = 39268.21
$ws.Range("I77").Value = 3961.8572
$ws.Range("J77").Value = 138126
$ws.Range("K77").Value = 19809.286
$ws.Range("L77").Value = 690630
$ws.Range("M77").Value = -15441.286
$ws.Range("N77").Value = -699366
$ws.Range("H102").Value = 3175.7368
$ws.Range("I102").Value = 2840.9333
$ws.Range("J102").Value = 4431.25
$ws.Range("K102").Value = 2840.9333
$ws.Range("L102").Value = 4431.25
$ws.Range("M102").Value = -1218.9333
$ws.Range("N102").Value = -7675.25
$ws.Range("H122").Value = 2242.389
$ws.Range("J122").Value = 4496.5
$ws.Range("L122").Value = 13489.5
$ws.Range("N122").Value = -18389.5

# ---- Worksheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1503
$ws.Range("I20").Value = 1503
$ws.Range("K20").Value = 1503
$ws.Range("M20").Value = -1256
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H107").Value = 4910.3335
$ws.Range("J107").Value = 3600
$ws.Range("L107").Value = 3600
$ws.Range("N107").Value = -7440
$ws.Range("H134").Value = 3485.348
$ws.Range("I134").Value = 3422.8865
$ws.Range("K134").Value = 10268.6595
$ws.Range("M134").Value = -7733.6595

# ---- Worksheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4121.1143
$ws.Range("I31").Value = 3006.6843
$ws.Range("K31").Value = 3006.6843
$ws.Range("M31").Value = -2711.6843
$ws.Range("H34").Value = 4121.1143
$ws.Range("I34").Value = 3006.6843
$ws.Range("K34").Value = 3006.6843
$ws.Range("M34").Value = -2804.6843
$ws.Range("H133").Value = 85000
$ws.Range("J133").Value = 85000
$ws.Range("L133").Value = 85000
$ws.Range("N133").Value = -90060
$ws.Range("H134").Value = 5765.6484
$ws.Range("I134").Value = 2933.8635
$ws.Range("K134").Value = 8801.5905
$ws.Range("M134").Value = -6266.5905

# ---- Worksheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 999.0476
$ws.Range("I5").Value = 759.4761999999999
$ws.Range("J5").Value = 1238.619
$ws.Range("K5").Value = 2278.4286
$ws.Range("L5").Value = 3715.857
$ws.Range("M5").Value = -2166.4286
$ws.Range("N5").Value = -3939.857
$ws.Range("H12").Value = 118.4
$ws.Range("I12").Value = 262.75
$ws.Range("J12").Value = 22.166666
$ws.Range("K12").Value = 788.25
$ws.Range("L12").Value = 66.49999800000001
$ws.Range("M12").Value = -615.25
$ws.Range("N12").Value = -412.499998
$ws.Range("H70").Value = 1001.25
$ws.Range("I70").Value = 1035
$ws.Range("K70").Value = 3105
$ws.Range("M70").Value = -2790
$ws.Range("H73").Value = 1001.25
$ws.Range("I73").Value = 1035
$ws.Range("K73").Value = 3105
$ws.Range("M73").Value = -2013
$ws.Range("H104").Value = 2804.3333
$ws.Range("J104").Value = 2750
$ws.Range("L104").Value = 8250
$ws.Range("N104").Value = -13492
$ws.Range("H135").Value = 999.0476
$ws.Range("I135").Value = 759.4761999999999
$ws.Range("J135").Value = 1238.619
$ws.Range("K135").Value = 6835.2858
$ws.Range("L135").Value = 11147.571
$ws.Range("M135").Value = -4300.2858
$ws.Range("N135").Value = -16217.571
$ws.Range("H138").Value = 1188.8
$ws.Range("I138").Value = 1188.8
$ws.Range("K138").Value = 3566.4
$ws.Range("M138").Value = 1573.6

# ---- Worksheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 169540.12
$ws.Range("I113").Value = 176781
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 176781
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -174611
$ws.Range("N113").Value = -7340
$ws.Range("H122").Value = 7442.2856
$ws.Range("I122").Value = 8579.4
$ws.Range("J122").Value = 4599.5
$ws.Range("K122").Value = 25738.2
$ws.Range("L122").Value = 13798.5
$ws.Range("M122").Value = -23288.2
$ws.Range("N122").Value = -18698.5
$ws.Range("H132").Value = 146
$ws.Range("J132").Value = 146
$ws.Range("L132").Value = 438
$ws.Range("N132").Value = -5498

# ---- Worksheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2365.3333
$ws.Range("I22").Value = 1147.5
$ws.Range("J22").Value = 2713.2856
$ws.Range("K22").Value = 1147.5
$ws.Range("L22").Value = 2713.2856
$ws.Range("M22").Value = -852.5
$ws.Range("N22").Value = -3303.2856
$ws.Range("H27").Value = 2365.3333
$ws.Range("I27").Value = 1147.5
$ws.Range("J27").Value = 2713.2856
$ws.Range("K27").Value = 1147.5
$ws.Range("L27").Value = 2713.2856
$ws.Range("M27").Value = -1040.5
$ws.Range("N27").Value = -2927.2856
$ws.Range("H46").Value = 3519.6
$ws.Range("I46").Value = 2349
$ws.Range("J46").Value = 3812.25
$ws.Range("K46").Value = 2349
$ws.Range("L46").Value = 3812.25
$ws.Range("M46").Value = -2161
$ws.Range("N46").Value = -4188.25
$ws.Range("H55").Value = 1890.7142
$ws.Range("I55").Value = 2138.9092
$ws.Range("J55").Value = 980.6667
$ws.Range("K55").Value = 2138.9092
$ws.Range("L55").Value = 980.6667
$ws.Range("M55").Value = -1965.9092
$ws.Range("N55").Value = -1326.6667
$ws.Range("H61").Value = 55138.35
$ws.Range("I61").Value = 60487.055
$ws.Range("K61").Value = 60487.055
$ws.Range("M61").Value = -60285.055
$ws.Range("H113").Value = 55138.35
$ws.Range("I113").Value = 60487.055
$ws.Range("K113").Value = 60487.055
$ws.Range("M113").Value = -58317.055

# ---- Worksheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 28372.25
$ws.Range("J45").Value = 24996.334
$ws.Range("L45").Value = 24996.334
$ws.Range("N45").Value = -25978.334
$ws.Range("H81").Value = 6491.25
$ws.Range("I81").Value = 2900
$ws.Range("J81").Value = 10082.5
$ws.Range("K81").Value = 5800
$ws.Range("L81").Value = 20165
$ws.Range("M81").Value = -4739
$ws.Range("N81").Value = -22287
$ws.Range("H84").Value = 6491.25
$ws.Range("I84").Value = 2900
$ws.Range("J84").Value = 10082.5
$ws.Range("K84").Value = 29000
$ws.Range("L84").Value = 100825
$ws.Range("M84").Value = -23696
$ws.Range("N84").Value = -111433
$ws.Range("H107").Value = 10676
$ws.Range("I107").Value = 10751
$ws.Range("K107").Value = 32253
$ws.Range("M107").Value = -30333
$ws.Range("H123").Value = 92000
$ws.Range("J123").Value = 92000
$ws.Range("L123").Value = 92000
$ws.Range("N123").Value = -101800

